# Add the new "ssh / hydra / fcrackzip / ncrack / hosts" bash-lib entries
# (rows 45-53) to 工作表1, matching the author's commit:
# "Add different bash entry and index entry that I learn when updating
#  elyn's website"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- row 45: Check ssh fail attempt on host -------------------------------
$ws.Range("A45").Value = "Linux"
$ws.Range("B45").Value = "Check ssh fail atttemp on host"
$ws.Range("C45").Value = "`$ grep sshd.*Failed /var/log/auth.log      //not only ssh login attemps, sudo attemp and other authoritzation related log "
$ws.Range("C45").WrapText = $true
$ws.Rows.Item(45).RowHeight = 32.4

# --- row 46: Adduser --------------------------------------------------------
$ws.Range("A46").Value = "Adduser"
$ws.Range("B46").Value = "Adduser to host"
$ws.Range("C46").Value = "`$ adduser new_user_name"
$ws.Range("C46").WrapText = $true

# --- row 47: ssh - update public key of a host ------------------------------
$ws.Range("A47").Value = "ssh"
$ws.Range("B47").Value = "Update public key of a host"
$ws.Range("C47").Value = "Update file: ~/.ssh/known_hosts`n* if a public key from a known host is updated, simply delete the old known_host entry"
$ws.Range("C47").WrapText = $true
$ws.Rows.Item(47).RowHeight = 32.4

# --- row 48: hydra -----------------------------------------------------------
$ws.Range("A48").Value = "hydra"
$ws.Range("B48").Value = "How to use hydra to crack a ssh"
$ws.Range("C48").Value = "1. Use nmap to find out whether ssh is on host and the possible username: `$nmap -sS -A ip_address`n2. Prepare a wordlist.txt (potential password) support the atk`n3. Use THC-hydra the tool to hack the target ssh: `$ hydra -user {username} -list {path_to_wordlist} {target_ip} ssh`nref: youtube/chris haralson: How to crack an SSH password"
$ws.Range("C48").WrapText = $true
$ws.Rows.Item(48).RowHeight = 97.2

# --- row 49: fcrackzip ---------------------------------------------------
$ws.Range("A49").Value = "fcrackzip"
$ws.Range("B49").Value = "Intro"
$ws.Range("C49").Value = "fcrackzip is an util to crack encrypted zip file. "
$ws.Range("C49").WrapText = $true

# --- row 50: ssh - basic ------------------------------------------------
$ws.Range("A50").Value = "ssh"
$ws.Range("B50").Value = "Basic"
$ws.Range("C50").Value = "`$ ssh -v myles@hostname //v for verbose"
$ws.Range("C50").WrapText = $true

# --- row 51: ssh - timeout troubleshooting -------------------------------
$ws.Range("A51").Value = "ssh"
$ws.Range("B51").Value = "Trouble-shoot ""Connection reset by xxx.xxx.xxx.xxx"" after a certain idle time at client side"
$ws.Range("C51").Value = "ref: https://www.bjornjohansen.no/ssh-timeout`nA. PREVENT SSH TIMEOUT ON THE CLIENT SIDE`nEdit your local SSH config file in ~/.ssh/config and add the following line:`n`` ServerAliveInterval 120`n(This will send a “null packet” every 120 seconds on your SSH connections to keep them alive.)`nB. PREVENT SSH TIMEOUT ON THE SERVER SIDE`nIf you’re a server admin, you can add the following to your SSH daemon config in /etc/ssh/sshd_config on your servers to prevent the clients to time out – so they don’t have to modify their local SSH config:`n`` ClientAliveInterval 120`n`` ClientAliveCountMax 720`nThis will make the server send the clients a “null packet” every 120 seconds and not disconnect them until the client have been inactive for 720 intervals (120 seconds * 720 = 86400 seconds = 24 hours)."
$ws.Range("C51").WrapText = $true
$ws.Rows.Item(51).RowHeight = 259.2

# --- row 52: ncrack --------------------------------------------------------
$ws.Range("A52").Value = "ncrack"
$ws.Range("B52").Value = "Intro"
$ws.Range("C52").Value = "Ncrack is an sibling project of nmap as another cracking tool like hydra, but it seems not working well when I try to apply it on my linode sshd (no login attemp when check auth.log on host)`n`$ ncrack --user myles {target_ip}:{port}`nsource code on github.com"
$ws.Range("C52").WrapText = $true
$ws.Rows.Item(52).RowHeight = 81

# --- row 53: linux - setup hosts ------------------------------------------
$ws.Range("A53").Value = "linux"
$ws.Range("B53").Value = "Setup Hosts and its benefic"
$ws.Range("C53").Value = "Setup hosts in file /etc/hosts. So that when in ping or ssh or other util, you can use the host alias directly like: `$ ssh myles@myleslinode"
$ws.Range("C53").WrapText = $true
$ws.Rows.Item(53).RowHeight = 32.4

# Widen column C to fit the new, longer content.
$ws.Columns.Item(3).ColumnWidth = 72

# Leave the cursor on the next empty row, same as the author did.
$ws.Range("C54").Select()
